$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C202").Value = 10402219140083218432.0
$ws.Range("C203").Value = 10388349514563106816.0
$ws.Range("C204").Value = 10374479889042995200.0
$ws.Range("C205").Value = 10360610263522885632.0
$ws.Range("C206").Value = 10346740638002774016.0
$ws.Range("C207").Value = 10332871012482662400.0
$ws.Range("C208").Value = 10319001386962552832.0
$ws.Range("C209").Value = 10305131761442441216.0
$ws.Range("C210").Value = 10291262135922329600.0
$ws.Range("C211").Value = 10277392510402217984.0
$ws.Range("C212").Value = 10263522884882108416.0
$ws.Range("C213").Value = 10249653259361998848.0
$ws.Range("C214").Value = 10235783633841887232.0
$ws.Range("C215").Value = 10221914008321773568.0
$ws.Range("C216").Value = 10208044382801664000.0
$ws.Range("C217").Value = 10194174757281554432.0
$ws.Range("C218").Value = 10180305131761442816.0
$ws.Range("C219").Value = 10166435506241331200.0
$ws.Range("C220").Value = 10152565880721219584.0
$ws.Range("C221").Value = 10138696255201110016.0
$ws.Range("C222").Value = 10124826629680998400.0
$ws.Range("C223").Value = 10110957004160886784.0
$ws.Range("C224").Value = 10097087378640777216.0
$ws.Range("C225").Value = 10083217753120665600.0
$ws.Range("C226").Value = 10069348127600553984.0
$ws.Range("C227").Value = 10055478502080442368.0
$ws.Range("C228").Value = 10041608876560332800.0
$ws.Range("C229").Value = 10027739251040223232.0
$ws.Range("C230").Value = 10013869625520111616.0
$ws.Range("C232").Value = 9986130374479888384.0
$ws.Range("C233").Value = 9972260748959778816.0
$ws.Range("C234").Value = 9958391123439667200.0
$ws.Range("C235").Value = 9944521497919555584.0
$ws.Range("C236").Value = 9930651872399446016.0
$ws.Range("C237").Value = 9916782246879334400.0
$ws.Range("C238").Value = 9902912621359222784.0
$ws.Range("C239").Value = 9889042995839111168.0
$ws.Range("C240").Value = 9875173370319001600.0
$ws.Range("C241").Value = 9861303744798892032.0
$ws.Range("C242").Value = 9847434119278780416.0
$ws.Range("C243").Value = 9833564493758668800.0
$ws.Range("C244").Value = 9819694868238557184.0
$ws.Range("C245").Value = 9805825242718447616.0
$ws.Range("C246").Value = 9791955617198336000.0
$ws.Range("C247").Value = 9778085991678224384.0
$ws.Range("C248").Value = 9764216366158112768.0
$ws.Range("C249").Value = 9750346740638003200.0
$ws.Range("C250").Value = 9736477115117891584.0
$ws.Range("C251").Value = 9722607489597779968.0
$ws.Range("C252").Value = 9708737864077670400.0
$ws.Range("C253").Value = 9694868238557558784.0
$ws.Range("C254").Value = 9680998613037447168.0
$ws.Range("C255").Value = 9667128987517335552.0
$ws.Range("C256").Value = 9653259361997225984.0
$ws.Range("C257").Value = 9639389736477116416.0
$ws.Range("C258").Value = 9625520110957004800.0
$ws.Range("C259").Value = 9611650485436893184.0
$ws.Range("C260").Value = 9597780859916781568.0
$ws.Range("C261").Value = 9583911234396672000.0
$ws.Range("C262").Value = 9570041608876560384.0
$ws.Range("C263").Value = 9556171983356448768.0
$ws.Range("C264").Value = 9542302357836339200.0
$ws.Range("C265").Value = 9528432732316227584.0
$ws.Range("C266").Value = 9514563106796115968.0
$ws.Range("C267").Value = 9500693481276004352.0
$ws.Range("C268").Value = 9486823855755894784.0
$ws.Range("C269").Value = 9472954230235785216.0
$ws.Range("C270").Value = 9459084604715673600.0
$ws.Range("C271").Value = 9445214979195561984.0
$ws.Range("C272").Value = 9431345353675450368.0
$ws.Range("C273").Value = 9417475728155340800.0
$ws.Range("C274").Value = 9403606102635229184.0
$ws.Range("C275").Value = 9389736477115117568.0
$ws.Range("C276").Value = 9375866851595005952.0
$ws.Range("C277").Value = 9361997226074896384.0
